# Auto-generated edit script: update cryptos.xlsx price/volume data per diff
# (hourly refresh of the cryptocurrency ranking table, incl. a few rows that
#  swapped position in the source ranking between runs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # The Price/Volume columns are plain text in the workbook (e.g. '1.00',
    # '0.526', multi-dot numbers like '69.362.09', or padded percentages).
    # Excel's COM Value setter auto-coerces simple numeric-looking strings
    # into real numbers (dropping formatting like trailing zeros), so we
    # briefly force Text number formatting for the assignment, then restore
    # the cell's original (default) style so no formatting diff is introduced.
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '69.362.09'
Set-TextValue "E2" '  -0.35%  '

# Row 3
Set-TextValue "D3" '3.741.27'
Set-TextValue "E3" '  -0.02%  '

# Row 4
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  -0.09%  '

# Row 5
Set-TextValue "D5" '615.15'
Set-TextValue "E5" '  +0.76%  '

# Row 6
Set-TextValue "D6" '178.57'
Set-TextValue "E6" '  +0.70%  '

# Row 7
Set-TextValue "D7" '3.740.77'
Set-TextValue "E7" '  -0.03%  '

# Row 8
Set-TextValue "E8" '  +0.06%  '

# Row 9
Set-TextValue "D9" '0.526'
Set-TextValue "E9" '  -1.32%  '

# Row 10
Set-TextValue "E10" '  +0.45%  '

# Row 11
Set-TextValue "D11" '6.54'
Set-TextValue "E11" '  +3.43%  '

# Row 12
Set-TextValue "D12" '0.481'
Set-TextValue "E12" '  -2.15%  '

# Row 13
Set-TextValue "D13" '39.97'
Set-TextValue "E13" '  -1.61%  '

# Row 14
Set-TextValue "D14" '0.0000254'
Set-TextValue "E14" '  +0.60%  '

# Row 15
Set-TextValue "D15" '4.362.64'

# Row 16
Set-TextValue "D16" '3.739.20'
Set-TextValue "E16" '  -0.18%  '

# Row 17
Set-TextValue "D17" '69.419.68'
Set-TextValue "E17" '  -0.35%  '

# Row 19
Set-TextValue "D19" '7.43'
Set-TextValue "E19" '  -1.54%  '

# Row 20
Set-TextValue "D20" '16.37'
Set-TextValue "E20" '  -1.50%  '

# Row 21
Set-TextValue "D21" '500.23'
Set-TextValue "E21" '  -2.35%  '

# Row 22
Set-TextValue "D22" '9.20'
Set-TextValue "E22" '  -2.98%  '

# Row 23
Set-TextValue "D23" '0.721'
Set-TextValue "E23" '  -0.20%  '

# Row 24
Set-TextValue "D24" '2.58'
Set-TextValue "E24" '  +3.30%  '

# Row 25
Set-TextValue "D25" '85.65'
Set-TextValue "E25" '  -2.32%  '

# Row 26
Set-TextValue "D26" '12.92'
Set-TextValue "E26" '  -2.35%  '

# Row 27
Set-TextValue "D27" '11.02'
Set-TextValue "E27" '  -0.15%  '

# Row 28
Set-TextValue "D28" '0.0000135'
Set-TextValue "E28" '  +6.20%  '

# Row 29
Set-TextValue "E29" '  -0.03%  '

# Row 30
Set-TextValue "D30" '2.50'
Set-TextValue "E30" '  +0.91%  '

# Row 31
Set-TextValue "B31" 'NEARProtocol'
Set-TextValue "C31" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D31" '8.07'
Set-TextValue "E31" '  +3.64%  '

# Row 32
Set-TextValue "B32" 'PancakeSwap'
Set-TextValue "C32" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D32" '2.91'
Set-TextValue "E32" '  +3.00%  '

# Row 33
Set-TextValue "D33" '30.39'
Set-TextValue "E33" '  -2.92%  '

# Row 34
Set-TextValue "E34" '  -1.29%  '

# Row 35
Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  +0.01%  '

# Row 36
Set-TextValue "B36" 'Filecoin'
Set-TextValue "C36" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D36" '6.12'
Set-TextValue "E36" '  -1.08%  '

# Row 37
Set-TextValue "B37" 'Mantle'
Set-TextValue "C37" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D37" '1.04'
Set-TextValue "E37" '  +0.32%  '

# Row 38
Set-TextValue "D38" '0.349'
Set-TextValue "E38" '  +3.63%  '

# Row 39
Set-TextValue "D39" '0.137'
Set-TextValue "E39" '  +3.89%  '

# Row 40
Set-TextValue "B40" 'Bittensor'
Set-TextValue "C40" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D40" '453.13'
Set-TextValue "E40" '  +8.50%  '

# Row 41
Set-TextValue "B41" 'dogwifhat'
Set-TextValue "C41" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D41" '3.08'
Set-TextValue "E41" '  +13.55%  '

# Row 42
Set-TextValue "D42" '2.07'
Set-TextValue "E42" '  -4.31%  '

# Row 43
Set-TextValue "D43" '49.69'
Set-TextValue "E43" '  -2.83%  '

# Row 44
Set-TextValue "D44" '44.72'
Set-TextValue "E44" '  +0.45%  '

# Row 45
Set-TextValue "D45" '8.56'
Set-TextValue "E45" '  -2.22%  '

# Row 46
Set-TextValue "D46" '2.947.94'
Set-TextValue "E46" '  -3.87%  '

# Row 47
Set-TextValue "D47" '0.0360'
Set-TextValue "E47" '  -0.58%  '

# Row 48
Set-TextValue "B48" 'Monero'
Set-TextValue "C48" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D48" '138.94'
Set-TextValue "E48" '  +2.30%  '

# Row 49
Set-TextValue "B49" 'InjectiveProtocol'
Set-TextValue "C49" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D49" '27.20'
Set-TextValue "E49" '  -1.56%  '

# Row 50
Set-TextValue "B50" 'USDe'
Set-TextValue "C50" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D50" '1.00'
Set-TextValue "E50" '  +0.00%  '

# Row 51
Set-TextValue "D51" '2.47'
Set-TextValue "E51" '  -1.62%  '
